$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Annotation scores for columns E:J (Clear, Assertive, Cautious, Optimistic,
# Specific, Relevant) for the 5 annotated rows (2-6).
$scores = @{
    2 = @(2, 2, 1, 2, 2, 2)
    3 = @(2, 2, 1, 2, 2, 2)
    4 = @(2, 2, 1, 1, 1, 2)
    5 = @(2, 2, 2, 2, 2, 2)
    6 = @(2, 2, 2, 2, 2, 2)
}

foreach ($row in $scores.Keys) {
    $vals = $scores[$row]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 5 + $j  # column E is index 5
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

# Move the window into view and update the active selection / zoom level
# to match the state the workbook was left in after annotating.
$excel.ActiveWindow.Left = -120
$ws.Range("G5").Select()
$excel.ActiveWindow.Zoom = 85
